$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # 'Citywide Totals'
$ws.Range("C2").Value = 55
$ws.Range("G2").Value = 69
$ws.Range("J2").Value = 90
$ws.Range("E3").Value = 107
$ws.Range("F3").Value = 100
$ws.Range("I3").Value = 165
$ws.Range("J3").Value = 172
$ws.Range("B9").Value = 301
$ws.Range("C9").Value = 370
$ws.Range("D9").Value = 327
$ws.Range("E9").Value = 346
$ws.Range("F9").Value = 409
$ws.Range("G9").Value = 380
$ws.Range("H9").Value = 358
$ws.Range("I9").Value = 409
$ws.Range("J9").Value = 322
$ws.Range("B10").Value = 1019
$ws.Range("C10").Value = 1230
$ws.Range("D10").Value = 1402
$ws.Range("E10").Value = 1708
$ws.Range("F10").Value = 1726
$ws.Range("G10").Value = 796
$ws.Range("H10").Value = 448
$ws.Range("I10").Value = 691
$ws.Range("J10").Value = 561
$ws.Range("K10").Value = 557
$ws.Range("B11").Value = 1429
$ws.Range("C11").Value = 1732
$ws.Range("D11").Value = 1917
$ws.Range("E11").Value = 2227
$ws.Range("F11").Value = 2311
$ws.Range("G11").Value = 1357
$ws.Range("H11").Value = 1017
$ws.Range("I11").Value = 1384
$ws.Range("J11").Value = 1169
$ws.Range("K11").Value = 1277

$ws = $wb.Worksheets.Item(11)  # 'Chicago Lawn'
$ws.Range("B6").Value = 9
$ws.Range("D6").Value = 9
$ws.Range("B7").Value = 22
$ws.Range("D7").Value = 19

$ws = $wb.Worksheets.Item(12)  # 'Garfield Park'
$ws.Range("J3").Value = 11
$ws.Range("D8").Value = 38
$ws.Range("H8").Value = 24
$ws.Range("D9").Value = 80
$ws.Range("H9").Value = 65
$ws.Range("J9").Value = 61

$ws = $wb.Worksheets.Item(14)  # 'Grand Crossing'
$ws.Range("B7").Value = 19
$ws.Range("B9").Value = 52

$ws = $wb.Worksheets.Item(15)  # 'Loop'
$ws.Range("I3").Value = 27
$ws.Range("B9").Value = 145
$ws.Range("C9").Value = 234
$ws.Range("D9").Value = 409
$ws.Range("E9").Value = 494
$ws.Range("I9").Value = 162
$ws.Range("K9").Value = 85
$ws.Range("B10").Value = 183
$ws.Range("C10").Value = 276
$ws.Range("D10").Value = 471
$ws.Range("E10").Value = 557
$ws.Range("I10").Value = 271
$ws.Range("K10").Value = 166

$ws = $wb.Worksheets.Item(16)  # 'Armour Square'
$ws.Range("E3").Value = 4
$ws.Range("I5").Value = 11
$ws.Range("F6").Value = 19
$ws.Range("E7").Value = 16
$ws.Range("F7").Value = 28
$ws.Range("I7").Value = 26

$ws = $wb.Worksheets.Item(17)  # 'Old Town'
$ws.Range("E6").Value = 44
$ws.Range("E7").Value = 55

$ws = $wb.Worksheets.Item(18)  # 'Little Italy, UIC'
$ws.Range("C6").Value = 11
$ws.Range("C7").Value = 18

$ws = $wb.Worksheets.Item(19)  # 'North Lawndale'
$ws.Range("C2").Value = 2
$ws.Range("C7").Value = 29
$ws.Range("C8").Value = 45

$ws = $wb.Worksheets.Item(2)  # 'By Neighborhood'
$ws.Range("E5").Value = 16
$ws.Range("F5").Value = 28
$ws.Range("I5").Value = 26
$ws.Range("H6").Value = 4
$ws.Range("D8").Value = 48
$ws.Range("E8").Value = 86
$ws.Range("J8").Value = 53
$ws.Range("C18").Value = 2
$ws.Range("B20").Value = 22
$ws.Range("D20").Value = 19
$ws.Range("B27").Value = 15
$ws.Range("B28").Value = 79
$ws.Range("D28").Value = 86
$ws.Range("E28").Value = 74
$ws.Range("F28").Value = 92
$ws.Range("D32").Value = 80
$ws.Range("H32").Value = 65
$ws.Range("J32").Value = 61
$ws.Range("B36").Value = 52
$ws.Range("G41").Value = 24
$ws.Range("J41").Value = 26
$ws.Range("F45").Value = 20
$ws.Range("C47").Value = 59
$ws.Range("J47").Value = 31
$ws.Range("C50").Value = 18
$ws.Range("C52").Value = 26
$ws.Range("F52").Value = 19
$ws.Range("B53").Value = 183
$ws.Range("C53").Value = 276
$ws.Range("D53").Value = 471
$ws.Range("E53").Value = 557
$ws.Range("I53").Value = 271
$ws.Range("K53").Value = 166
$ws.Range("K54").Value = 8
$ws.Range("F56").Value = 11
$ws.Range("D61").Value = 20
$ws.Range("C65").Value = 45
$ws.Range("E70").Value = 55
$ws.Range("J74").Value = 31
$ws.Range("C76").Value = 60
$ws.Range("H77").Value = 44
$ws.Range("B78").Value = 33
$ws.Range("E78").Value = 41
$ws.Range("D87").Value = 18
$ws.Range("G88").Value = 10
$ws.Range("G95").Value = 14
$ws.Range("F98").Value = 10
$ws.Range("B99").Value = 1429
$ws.Range("C99").Value = 1732
$ws.Range("D99").Value = 1917
$ws.Range("E99").Value = 2227
$ws.Range("F99").Value = 2311
$ws.Range("G99").Value = 1357
$ws.Range("H99").Value = 1017
$ws.Range("I99").Value = 1384
$ws.Range("J99").Value = 1169
$ws.Range("K99").Value = 1277

$ws = $wb.Worksheets.Item(22)  # 'Humboldt Park'
$ws.Range("J5").Value = 9
$ws.Range("G6").Value = 18
$ws.Range("G7").Value = 24
$ws.Range("J7").Value = 26

$ws = $wb.Worksheets.Item(24)  # 'Uptown'
$ws.Range("D8").Value = 15
$ws.Range("D9").Value = 18

$ws = $wb.Worksheets.Item(25)  # 'Rush & Division'
$ws.Range("E3").Value = 2
$ws.Range("B5").Value = 28
$ws.Range("B6").Value = 33
$ws.Range("E6").Value = 41

$ws = $wb.Worksheets.Item(26)  # 'Englewood'
$ws.Range("B7").Value = 27
$ws.Range("E7").Value = 18
$ws.Range("D8").Value = 49
$ws.Range("F8").Value = 51
$ws.Range("B9").Value = 79
$ws.Range("D9").Value = 86
$ws.Range("E9").Value = 74
$ws.Range("F9").Value = 92

$ws = $wb.Worksheets.Item(27)  # 'Lake View'
$ws.Range("J2").Value = 2
$ws.Range("C6").Value = 12
$ws.Range("C8").Value = 59
$ws.Range("J8").Value = 31

$ws = $wb.Worksheets.Item(28)  # 'Jefferson Park'
$ws.Range("F6").Value = 17
$ws.Range("F7").Value = 20

$ws = $wb.Worksheets.Item(3)  # 'Rogers Park'
$ws.Range("C7").Value = 48
$ws.Range("C8").Value = 60

$ws = $wb.Worksheets.Item(31)  # 'River North'
$ws.Range("J6").Value = 11
$ws.Range("J7").Value = 31

$ws = $wb.Worksheets.Item(38)  # 'Logan Square'
$ws.Range("C7").Value = 19
$ws.Range("F7").Value = 15
$ws.Range("C8").Value = 26
$ws.Range("F8").Value = 19

$ws = $wb.Worksheets.Item(4)  # 'Edgewater'
$ws.Range("B5").Value = 2
$ws.Range("B7").Value = 15

$ws = $wb.Worksheets.Item(43)  # 'Lower West Side'
$ws.Range("K5").Value = 5
$ws.Range("K6").Value = 8

$ws = $wb.Worksheets.Item(45)  # 'Calumet Heights'
$ws.Range("C4").Value = 2
$ws.Range("C6").Value = 2

$ws = $wb.Worksheets.Item(5)  # 'Roseland'
$ws.Range("H8").Value = 13
$ws.Range("H10").Value = 44

$ws = $wb.Worksheets.Item(55)  # 'West Town'
$ws.Range("G2").Value = 1
$ws.Range("G7").Value = 14

$ws = $wb.Worksheets.Item(65)  # 'Wrigleyville'
$ws.Range("F3").Value = 1
$ws.Range("F6").Value = 4
$ws.Range("F7").Value = 10

$ws = $wb.Worksheets.Item(78)  # 'Ashburn'
$ws.Range("H4").Value = 1
$ws.Range("H6").Value = 4

$ws = $wb.Worksheets.Item(8)  # 'Austin'
$ws.Range("D6").Value = 18
$ws.Range("J6").Value = 18
$ws.Range("E7").Value = 47
$ws.Range("D8").Value = 48
$ws.Range("E8").Value = 86
$ws.Range("J8").Value = 53

$ws = $wb.Worksheets.Item(86)  # 'Mckinley Park'
$ws.Range("F4").Value = 1
$ws.Range("F6").Value = 11

$ws = $wb.Worksheets.Item(9)  # 'Washington Heights'
$ws.Range("G5").Value = 4
$ws.Range("G7").Value = 10
